$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95, shifting existing rows 95-100 down to 96-101.
$ws.Rows.Item(95).Insert()

# Copy the date number format (style) from the row above into the new row's D cell
$ws.Cells.Item(94, 4).Copy()
$ws.Cells.Item(95, 4).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new row 95 with the latest weekly entry
$ws.Cells.Item(95, 1).Value = 4
$ws.Cells.Item(95, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(95, 3).Value = "Los Lagos"
$ws.Cells.Item(95, 4).Value = 45265
$ws.Cells.Item(95, 5).Value = 10
$ws.Cells.Item(95, 6).Value = 100112012
$ws.Cells.Item(95, 7).Value = "Espinaca"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 40
$ws.Cells.Item(95, 11).Value = 20000
$ws.Cells.Item(95, 12).Value = 20000
$ws.Cells.Item(95, 13).Value = 20000
$ws.Cells.Item(95, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(95, 15).Value = "Región Metropolitana"
$ws.Cells.Item(95, 16).Value = 2000
$ws.Cells.Item(95, 17).Value = 10
$ws.Cells.Item(95, 18).Value = "Hortaliza"
